$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.027278619896107
$ws.Cells.Item(2, 4).Value = 1.032755763016446
$ws.Cells.Item(2, 5).Value = 1.051046160158886
$ws.Cells.Item(2, 6).Value = 1.056080361991812
$ws.Cells.Item(2, 9).Value = 1.035352699529619
$ws.Cells.Item(2, 10).Value = 1.032437380866821
$ws.Cells.Item(2, 11).Value = 1.035560077793327
$ws.Cells.Item(2, 12).Value = 1.053798749314762
$ws.Cells.Item(2, 13).Value = 1.058819055084351
$ws.Cells.Item(2, 14).Value = 1.014845711269824

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.028019124234036
$ws.Cells.Item(3, 4).Value = 1.033307649016017
$ws.Cells.Item(3, 5).Value = 1.052059733081065
$ws.Cells.Item(3, 6).Value = 1.057098913616316
$ws.Cells.Item(3, 9).Value = 1.035506303731186
$ws.Cells.Item(3, 10).Value = 1.032819076604272
$ws.Cells.Item(3, 11).Value = 1.035921489365779
$ws.Cells.Item(3, 12).Value = 1.054624372853338
$ws.Cells.Item(3, 13).Value = 1.059650652848738
$ws.Cells.Item(3, 14).Value = 1.014972518812272

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.028498789268428
$ws.Cells.Item(4, 4).Value = 1.033665159554896
$ws.Cells.Item(4, 5).Value = 1.052716905036241
$ws.Cells.Item(4, 6).Value = 1.057759166328671
$ws.Cells.Item(4, 9).Value = 1.035604773254044
$ws.Cells.Item(4, 10).Value = 1.03306587613926
$ws.Cells.Item(4, 11).Value = 1.036155043200257
$ws.Cells.Item(4, 12).Value = 1.055159316716236
$ws.Cells.Item(4, 13).Value = 1.060189321729036
$ws.Cells.Item(4, 14).Value = 1.01505449409862

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028700560604521
$ws.Cells.Item(5, 4).Value = 1.033815551807648
$ws.Cells.Item(5, 5).Value = 1.052993495480349
$ws.Cells.Item(5, 6).Value = 1.058037018122124
$ws.Cells.Item(5, 9).Value = 1.035645948291513
$ws.Cells.Item(5, 10).Value = 1.0331695856929
$ws.Cells.Item(5, 11).Value = 1.036253155361807
$ws.Cells.Item(5, 12).Value = 1.055384376030869
$ws.Cells.Item(5, 13).Value = 1.060415913249575
$ws.Cells.Item(5, 14).Value = 1.015088937511054

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028734445894792
$ws.Cells.Item(6, 4).Value = 1.033840808838081
$ws.Cells.Item(6, 5).Value = 1.053039954689152
$ws.Cells.Item(6, 6).Value = 1.05808368713219
$ws.Cells.Item(6, 9).Value = 1.035652848750483
$ws.Cells.Item(6, 10).Value = 1.033186996313354
$ws.Cells.Item(6, 11).Value = 1.036269624459229
$ws.Cells.Item(6, 12).Value = 1.055422174353054
$ws.Cells.Item(6, 13).Value = 1.060453966851245
$ws.Cells.Item(6, 14).Value = 1.015094719585789

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.02850148487696
$ws.Cells.Item(7, 4).Value = 1.033667168733158
$ws.Cells.Item(7, 5).Value = 1.052720599613728
$ws.Cells.Item(7, 6).Value = 1.057762877893301
$ws.Cells.Item(7, 9).Value = 1.035605324308273
$ws.Cells.Item(7, 10).Value = 1.033067262089513
$ws.Cells.Item(7, 11).Value = 1.036156354471507
$ws.Cells.Item(7, 12).Value = 1.055162323307482
$ws.Cells.Item(7, 13).Value = 1.060192348926603
$ws.Cells.Item(7, 14).Value = 1.015054954408338

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.027528770315698
$ws.Cells.Item(8, 4).Value = 1.032942190896851
$ws.Cells.Item(8, 5).Value = 1.051388427285022
$ws.Cells.Item(8, 6).Value = 1.056424340789918
$ws.Cells.Item(8, 9).Value = 1.03540480143475
$ws.Cells.Item(8, 10).Value = 1.032566413968019
$ws.Cells.Item(8, 11).Value = 1.035682280681658
$ws.Cells.Item(8, 12).Value = 1.054077624926357
$ws.Cells.Item(8, 13).Value = 1.059099978783919
$ws.Cells.Item(8, 14).Value = 1.014888582313116

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.025818700383321
$ws.Cells.Item(9, 4).Value = 1.031667854170087
$ws.Cells.Item(9, 5).Value = 1.049051153833186
$ws.Cells.Item(9, 6).Value = 1.054074777626129
$ws.Cells.Item(9, 9).Value = 1.035044417294823
$ws.Cells.Item(9, 10).Value = 1.031682502941437
$ws.Cells.Item(9, 11).Value = 1.034844629030432
$ws.Cells.Item(9, 12).Value = 1.052171736700644
$ws.Cells.Item(9, 13).Value = 1.057179498987826
$ws.Cells.Item(9, 14).Value = 1.014594836341686

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.024681432385508
$ws.Cells.Item(10, 4).Value = 1.030820530049268
$ws.Cells.Item(10, 5).Value = 1.047499893311798
$ws.Cells.Item(10, 6).Value = 1.05251460740559
$ws.Cells.Item(10, 9).Value = 1.03479946931918
$ws.Cells.Item(10, 10).Value = 1.031092386807783
$ws.Cells.Item(10, 11).Value = 1.034284736906163
$ws.Cells.Item(10, 12).Value = 1.050904901093892
$ws.Cells.Item(10, 13).Value = 1.055902220654183
$ws.Cells.Item(10, 14).Value = 1.014398641705447

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.024189663873604
$ws.Cells.Item(11, 4).Value = 1.030454181272207
$ws.Cells.Item(11, 5).Value = 1.046829837602976
$ws.Cells.Item(11, 6).Value = 1.051840524942408
$ws.Cells.Item(11, 9).Value = 1.034692299623551
$ws.Cells.Item(11, 10).Value = 1.030836674473837
$ws.Cells.Item(11, 11).Value = 1.034041966598046
$ws.Cells.Item(11, 12).Value = 1.050357252052732
$ws.Cells.Item(11, 13).Value = 1.055349881562941
$ws.Cells.Item(11, 14).Value = 1.014313605995539

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.024007102496063
$ws.Cells.Item(12, 4).Value = 1.030318187176243
$ws.Cells.Item(12, 5).Value = 1.046581198246616
$ws.Cells.Item(12, 6).Value = 1.051590364606902
$ws.Cells.Item(12, 9).Value = 1.034652326580105
$ws.Cells.Item(12, 10).Value = 1.030741664400942
$ws.Cells.Item(12, 11).Value = 1.033951742017045
$ws.Cells.Item(12, 12).Value = 1.050153966984972
$ws.Cells.Item(12, 13).Value = 1.05514482922883
$ws.Cells.Item(12, 14).Value = 1.014282008013886

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024046257826204
$ws.Cells.Item(13, 4).Value = 1.030347354534632
$ws.Cells.Item(13, 5).Value = 1.046634520937473
$ws.Cells.Item(13, 6).Value = 1.051644014694619
$ws.Cells.Item(13, 9).Value = 1.03466090840985
$ws.Cells.Item(13, 10).Value = 1.030762045593306
$ws.Cells.Item(13, 11).Value = 1.033971097688594
$ws.Cells.Item(13, 12).Value = 1.0501975661426
$ws.Cells.Item(13, 13).Value = 1.055188808608463
$ws.Cells.Item(13, 14).Value = 1.014288786421659

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.024174571174724
$ws.Cells.Item(14, 4).Value = 1.030442938230769
$ws.Cells.Item(14, 5).Value = 1.046809279906157
$ws.Cells.Item(14, 6).Value = 1.051819842035784
$ws.Cells.Item(14, 9).Value = 1.034688998810562
$ws.Cells.Item(14, 10).Value = 1.030828821458345
$ws.Cells.Item(14, 11).Value = 1.034034509596226
$ws.Cells.Item(14, 12).Value = 1.050340445661487
$ws.Cells.Item(14, 13).Value = 1.055332929601796
$ws.Cells.Item(14, 14).Value = 1.014310994337589

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.024253643010843
$ws.Cells.Item(15, 4).Value = 1.030501841705185
$ws.Cells.Item(15, 5).Value = 1.046916987728469
$ws.Cells.Item(15, 6).Value = 1.051928204772895
$ws.Cells.Item(15, 9).Value = 1.034706284329949
$ws.Cells.Item(15, 10).Value = 1.030869960712293
$ws.Cells.Item(15, 11).Value = 1.034073573324247
$ws.Cells.Item(15, 12).Value = 1.050428496525132
$ws.Cells.Item(15, 13).Value = 1.055421742046433
$ws.Cells.Item(15, 14).Value = 1.014324675798518

# Row 16
$ws.Cells.Item(16, 2).Value = 1.019999999999999
$ws.Cells.Item(16, 3).Value = 1.024714083837009
$ws.Cells.Item(16, 4).Value = 1.030844855113611
$ws.Cells.Item(16, 5).Value = 1.047544397577384
$ws.Cells.Item(16, 6).Value = 1.052559375365495
$ws.Cells.Item(16, 9).Value = 1.034806558581212
$ws.Cells.Item(16, 10).Value = 1.031109353726482
$ws.Cells.Item(16, 11).Value = 1.034300841851582
$ws.Cells.Item(16, 12).Value = 1.050941265793368
$ws.Cells.Item(16, 13).Value = 1.055938893080693
$ws.Cells.Item(16, 14).Value = 1.014404283549702

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.025003088450811
$ws.Cells.Item(17, 4).Value = 1.031060166389875
$ws.Cells.Item(17, 5).Value = 1.047938397902904
$ws.Cells.Item(17, 6).Value = 1.052955689582777
$ws.Cells.Item(17, 9).Value = 1.034869162431792
$ws.Cells.Item(17, 10).Value = 1.031259469071463
$ws.Cells.Item(17, 11).Value = 1.034443313000886
$ws.Cells.Item(17, 12).Value = 1.051263153896624
$ws.Cells.Item(17, 13).Value = 1.056263484773741
$ws.Cells.Item(17, 14).Value = 1.014454197686028

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.025171725060934
$ws.Cells.Item(18, 4).Value = 1.031185806583639
$ws.Cells.Item(18, 5).Value = 1.048168370777619
$ws.Cells.Item(18, 6).Value = 1.053186995818593
$ws.Cells.Item(18, 9).Value = 1.034905571479529
$ws.Cells.Item(18, 10).Value = 1.031347010581795
$ws.Cells.Item(18, 11).Value = 1.034526381770617
$ws.Cells.Item(18, 12).Value = 1.051450992461846
$ws.Cells.Item(18, 13).Value = 1.056452884053198
$ws.Cells.Item(18, 14).Value = 1.014483303804867

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.025229236776497
$ws.Cells.Item(19, 4).Value = 1.031228655538705
$ws.Cells.Item(19, 5).Value = 1.048246812542461
$ws.Cells.Item(19, 6).Value = 1.053265889458638
$ws.Cells.Item(19, 9).Value = 1.034917967902499
$ws.Cells.Item(19, 10).Value = 1.031376856847571
$ws.Cells.Item(19, 11).Value = 1.034554700569244
$ws.Cells.Item(19, 12).Value = 1.051515055228002
$ws.Cells.Item(19, 13).Value = 1.056517476223994
$ws.Cells.Item(19, 14).Value = 1.014493226877953

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.02497207424357
$ws.Cells.Item(20, 4).Value = 1.031037060056498
$ws.Cells.Item(20, 5).Value = 1.047896108926192
$ws.Cells.Item(20, 6).Value = 1.052913154031925
$ws.Cells.Item(20, 9).Value = 1.034862456664819
$ws.Cells.Item(20, 10).Value = 1.031243364993069
$ws.Cells.Item(20, 11).Value = 1.034428030521112
$ws.Cells.Item(20, 12).Value = 1.051228609375845
$ws.Cells.Item(20, 13).Value = 1.056228651862765
$ws.Cells.Item(20, 14).Value = 1.014448843191196

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024136783224643
$ws.Cells.Item(21, 4).Value = 1.030414788877215
$ws.Cells.Item(21, 5).Value = 1.046757810857508
$ws.Cells.Item(21, 6).Value = 1.051768059074684
$ws.Cells.Item(21, 9).Value = 1.034680731449536
$ws.Cells.Item(21, 10).Value = 1.030809158370095
$ws.Cells.Item(21, 11).Value = 1.034015837709612
$ws.Cells.Item(21, 12).Value = 1.050298367436644
$ws.Cells.Item(21, 13).Value = 1.05529048648769
$ws.Cells.Item(21, 14).Value = 1.014304454985999

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.023612201403835
$ws.Cells.Item(22, 4).Value = 1.030024030049979
$ws.Cells.Item(22, 5).Value = 1.046043559903747
$ws.Cells.Item(22, 6).Value = 1.051049388529436
$ws.Cells.Item(22, 9).Value = 1.03456551669692
$ws.Cells.Item(22, 10).Value = 1.030535999051859
$ws.Cells.Item(22, 11).Value = 1.033756393521747
$ws.Cells.Item(22, 12).Value = 1.049714275788716
$ws.Cells.Item(22, 13).Value = 1.05470126747242
$ws.Cells.Item(22, 14).Value = 1.014213603552592

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.023890234799997
$ws.Cells.Item(23, 4).Value = 1.030231131864711
$ws.Cells.Item(23, 5).Value = 1.046422060731963
$ws.Cells.Item(23, 6).Value = 1.051430246084831
$ws.Cells.Item(23, 9).Value = 1.03462668469989
$ws.Cells.Item(23, 10).Value = 1.030680820416858
$ws.Cells.Item(23, 11).Value = 1.033893956144527
$ws.Cells.Item(23, 12).Value = 1.050023838729352
$ws.Cells.Item(23, 13).Value = 1.055013562291143
$ws.Cells.Item(23, 14).Value = 1.014261772005528

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.024986088027883
$ws.Cells.Item(24, 4).Value = 1.031047500651141
$ws.Cells.Item(24, 5).Value = 1.047915217002825
$ws.Cells.Item(24, 6).Value = 1.052932373575153
$ws.Cells.Item(24, 9).Value = 1.034865487042094
$ws.Cells.Item(24, 10).Value = 1.031250641789378
$ws.Cells.Item(24, 11).Value = 1.034434936115669
$ws.Cells.Item(24, 12).Value = 1.051244218290105
$ws.Cells.Item(24, 13).Value = 1.056244391138347
$ws.Cells.Item(24, 14).Value = 1.014451262681552

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.026260312349475
$ws.Cells.Item(25, 4).Value = 1.031996914691159
$ws.Cells.Item(25, 5).Value = 1.049654180503201
$ws.Cells.Item(25, 6).Value = 1.05468110819317
$ws.Cells.Item(25, 9).Value = 1.035138415110817
$ws.Cells.Item(25, 10).Value = 1.031911168485585
$ws.Cells.Item(25, 11).Value = 1.035061444153627
$ws.Cells.Item(25, 12).Value = 1.052663797236854
$ws.Cells.Item(25, 13).Value = 1.057675458287376
$ws.Cells.Item(25, 14).Value = 1.014670842628012

